# Actualizacion automatica del mapa (2025-10-30 07:13:40)
# Adds the new reclamo row (Caso 7690) at the bottom of the "NEW" sheet,
# mirroring the layout/types used by the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 81

# Force the text-like columns to be stored as text (not auto-coerced to
# numbers) before assigning their values, matching how the other rows in
# this sheet store "numeric-looking" identifiers (Caso, Comuna, OT, etc.)
# as strings rather than numbers.
$ws.Range("A" + $newRow + ":H" + $newRow).NumberFormat = "@"
$ws.Range("J" + $newRow + ":L" + $newRow).NumberFormat = "@"
$ws.Range("O" + $newRow + ":R" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "7690"
$ws.Range("B" + $newRow).Value = "10/29/2025"
$ws.Range("C" + $newRow).Value = "ESCALADA DE SAN MARTIN, R. 2555"
$ws.Range("D" + $newRow).Value = "11"
$ws.Range("E" + $newRow).Value = "810481198"
$ws.Range("F" + $newRow).Value = "NEW"
$ws.Range("G" + $newRow).Value = "Pendiente"
$ws.Range("H" + $newRow).Value = "Picada"
$ws.Range("I" + $newRow).Value = 1
$ws.Range("J" + $newRow).Value = "Cambio"
$ws.Range("K" + $newRow).Value = "Sin equipos"
$ws.Range("L" + $newRow).Value = "Pasante"
$ws.Range("M" + $newRow).Value = -58.473913
$ws.Range("N" + $newRow).Value = -34.613111
$ws.Range("O" + $newRow).Value = "Paternal"
$ws.Range("P" + $newRow).Value = "Capital Norte"
$ws.Range("Q" + $newRow).Value = "NRA-K"
$ws.Range("R" + $newRow).Value = "Fuera de Poligono OVL"

# Drop the temporary "Text" number format again so the new cells end up
# with the same (default) style as every other data row in the sheet.
$ws.Range("A" + $newRow + ":R" + $newRow).ClearFormats()
